# Updates cryptos list data (Coin / Link / Price / Volume(1h)) per
# the Fri Jun 14 03:48:54 UTC 2024 GitHub Actions refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.749.81"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "3.495.17"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'600.74"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "'147.63"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("D7").Value = "3.494.89"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'7.91"
$ws.Range("E11").Value = "  +5.68%  "
$ws.Range("D12").Value = "'0.421"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "4.086.00"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  -4.04%  "
$ws.Range("D16").Value = "3.490.17"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "66.815.63"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "'10.46"
$ws.Range("E19").Value = "  +7.03%  "
$ws.Range("D20").Value = "'6.36"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("D21").Value = "'15.29"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").Value = "'433.23"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("D24").Value = "'79.63"
$ws.Range("E24").Value = "  +3.00%  "
$ws.Range("D25").Value = "3.633.41"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("E28").Value = "  -6.48%  "
$ws.Range("D29").Value = "'9.82"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "'8.21"
$ws.Range("E30").Value = "  -7.01%  "
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").Value = "'1.60"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("D35").Value = "'25.32"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").Value = "3.488.93"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "'0.0891"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "'170.11"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("E44").Value = "  -8.99%  "
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "'0.896"
$ws.Range("E46").Value = "  +2.36%  "
$ws.Range("D47").Value = "'45.83"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("D48").Value = "'28.20"
$ws.Range("E48").Value = "  -6.12%  "
$ws.Range("D49").Value = "'1.29"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("E51").Value = "  -4.01%  "

# Rows 39/40: USDe and Aptos swap ranking positions, each with updated
# price/volume figures.
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'7.98"
$ws.Range("E40").Value = "  -0.23%  "
